# Update "想去人数" (interested-count) figures across the four sheets to the
# newly scraped values (gh-pages data refresh @ 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 118
$ws1.Range("F6").Value  = 614
$ws1.Range("F7").Value  = 329
$ws1.Range("F8").Value  = 528
$ws1.Range("F10").Value = 10813
$ws1.Range("F14").Value = 2043
$ws1.Range("F18").Value = 211
$ws1.Range("F19").Value = 60
$ws1.Range("F20").Value = 233
$ws1.Range("F21").Value = 1156
$ws1.Range("F22").Value = 128
$ws1.Range("F23").Value = 210
$ws1.Range("F24").Value = 694
$ws1.Range("F26").Value = 223
$ws1.Range("F27").Value = 2379
$ws1.Range("F28").Value = 676
$ws1.Range("F29").Value = 3197
$ws1.Range("F30").Value = 1015
$ws1.Range("F31").Value = 735
$ws1.Range("F35").Value = 930
$ws1.Range("F36").Value = 24
$ws1.Range("F37").Value = 29
$ws1.Range("F38").Value = 234
$ws1.Range("F39").Value = 2
$ws1.Range("F41").Value = 1275
$ws1.Range("F43").Value = 96
$ws1.Range("F44").Value = 132
$ws1.Range("F45").Value = 226
$ws1.Range("F47").Value = 11
$ws1.Range("F49").Value = 82

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 4099
$ws2.Range("F6").Value  = 78
$ws2.Range("F12").Value = 377

# --- 本地生活 (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 739
$ws3.Range("F3").Value = 412
$ws3.Range("F4").Value = 50

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 739
$ws4.Range("F3").Value  = 412
$ws4.Range("F4").Value  = 50
$ws4.Range("F6").Value  = 118
$ws4.Range("F7").Value  = 614
$ws4.Range("F9").Value  = 10813
$ws4.Range("F14").Value = 211
$ws4.Range("F15").Value = 233
$ws4.Range("F16").Value = 1156
$ws4.Range("F17").Value = 128
$ws4.Range("F18").Value = 210
$ws4.Range("F19").Value = 4099
$ws4.Range("F21").Value = 694
$ws4.Range("F23").Value = 223
$ws4.Range("F24").Value = 676
$ws4.Range("F25").Value = 3197
$ws4.Range("F26").Value = 1015
$ws4.Range("F27").Value = 78
$ws4.Range("F29").Value = 735
$ws4.Range("F32").Value = 24
$ws4.Range("F33").Value = 29
$ws4.Range("F34").Value = 234
$ws4.Range("F35").Value = 1275
$ws4.Range("F37").Value = 96
$ws4.Range("F38").Value = 132
$ws4.Range("F39").Value = 226
$ws4.Range("F43").Value = 11
$ws4.Range("F49").Value = 82
